# The BOM title (merged cell A1:D1 on Sheet1) is renamed from the old
# "KickJr" project title to the new "NotSoSmartWatch" project title/author,
# matching the repository being copied from KickJr -> NotSoSmartWatch.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Click on / select the merged title cell, like a user would before typing
# a replacement value.
$ws.Range("A1:D1").Select() | Out-Null

# Leading apostrophe forces explicit-text entry (keeps the cell's existing
# "quote prefix" text formatting) exactly like typing '... into Excel.
$ws.Range("A1").Value = "'NotSoSmartWatch PPG Module -- Rev A / Thomas Davis"
